$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Methods")

# The "Users" section gained ten more methods (rows 61-70, column B already
# held the new method paths). Mark each newly-added row as "Done" in column
# A, matching the existing pattern used throughout the checklist (shared
# "Done" text + the built-in "Good" cell style).
for ($r = 61; $r -le 70; $r++) {
    $cell = $ws.Range("A$r")
    $cell.Value = "Done"
    $cell.Style = "Good"
}

# Bring the newly completed rows into view / selection, same as the author
# scrolling down to the bottom of the freshly extended list.
$ws.Activate()
$ws.Range("A70").Select()
